$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.702.85"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.925.88"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'604.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'165.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "3.920.72"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'37.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "4.588.01"
$ws.Range("E15").Value = "  +3.75%  "
$ws.Range("D16").Value = "3.908.49"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "68.896.31"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'7.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'17.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'11.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'486.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "'0.0000171"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.22%  "
$ws.Range("D24").Value = "'0.724"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'12.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'10.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").Value = "4.080.05"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "'2.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'7.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'32.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").Value = "3.877.92"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'3.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'0.318"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").Value = "'434.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").Value = "'48.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D48").Value = "'26.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.42%  "
$ws.Range("D49").Value = "2.837.88"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("D50").Value = "'141.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'39.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
